$d = $word.ActiveDocument

# 1. Title heading + bold byline at the end (same text appears twice)
$d.Content.Find.Execute("Play Banana Splash Free - Review of Online Slot Game", $true, $false, $false, $false, $false, $true, 1, $false, "Play Banana Splash for Free - Fun and Flexible Beach-Themed Slot", 2)
$d.Content.Find.Execute("Play Banana Splash Free - Review of Online Slot Game", $true, $false, $false, $false, $false, $true, 1, $false, "Play Banana Splash for Free - Fun and Flexible Beach-Themed Slot", 2)

# 2. "What we like" bullets
$d.Content.Find.Execute("Simple and easy-to-play game", $true, $false, $false, $false, $false, $true, 1, $false, "Simple and easy to play for beginners", 2)
$d.Content.Find.Execute("Fun and colorful beachy graphics", $true, $false, $false, $false, $false, $true, 1, $false, "Fun and colorful beach theme with tasteful graphics", 2)
$d.Content.Find.Execute("Nine paylines available for players", $true, $false, $false, $false, $false, $true, 1, $false, "Flexible paylines and adjustable betting options", 2)
$d.Content.Find.Execute("Low minimum bet per spin", $true, $false, $false, $false, $false, $true, 1, $false, "Acceptable RTP and bonus features", 2)

# 3. "What we don't like" bullet
$d.Content.Find.Execute("Slightly below average return to player percentage", $true, $false, $false, $false, $false, $true, 1, $false, "Slightly below-average RTP", 2)

# 4. Final italic summary paragraph
$d.Content.Find.Execute("Learn more about Banana Splash - a fun and simple online slot game with colorful graphics. Play Banana Splash for free and enjoy its bonus rounds.", $true, $false, $false, $false, $false, $true, 1, $false, "Beginner-friendly online slot game with adjustable paylines, fun fruit symbols, and pleasing graphics.", 2)
